{"js": "const replacements = [\n  [\"637\u00f75=\", \"653\u00f73=\"],\n  [\"337\u00f79=\", \"902\u00f76=\"],\n  [\"379\u00f79=\", \"796\u00f79=\"],\n  [\"692\u00f79=\", \"595\u00f77=\"],\n  [\"913\u00f74=\", \"899\u00f78=\"],\n  [\"628\u00f78=\", \"103\u00f78=\"],\n  [\"975\u00f79=\", \"983\u00f73=\"],\n  [\"116\u00f77=\", \"972\u00f77=\"],\n  [\"386\u00f76=\", \"957\u00f76=\"],\n  [\"804\u00f78=\", \"893\u00f78=\"],\n  [\"998\u00f75=\", \"106\u00f75=\"],\n  [\"404\u00f78=\", \"628\u00f75=\"],\n  [\"836\u00f72=\", \"863\u00f74=\"],\n  [\"653\u00f74=\", \"665\u00f76=\"],\n  [\"637\u00f79=\", \"597\u00f74=\"],\n  [\"300\u00f76=\", \"763\u00f77=\"],\n  [\"843\u00f78=\", \"317\u00f73=\"],\n  [\"430\u00f74=\", \"418\u00f77=\"],\n  [\"211\u00f73=\", \"585\u00f77=\"],\n  [\"323\u00f75=\", \"989\u00f78=\"],\n  [\"676\u00f72=\", \"187\u00f73=\"],\n  [\"161\u00f74=\", \"999\u00f74=\"],\n  [\"996\u00f72=\", \"915\u00f75=\"],\n  [\"412\u00f74=\", \"453\u00f75=\"],\n  [\"240\u00f72=\", \"497\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"637\u00f75=\", \"653\u00f73=\"),\n    @(\"337\u00f79=\", \"902\u00f76=\"),\n    @(\"379\u00f79=\", \"796\u00f79=\"),\n    @(\"692\u00f79=\", \"595\u00f77=\"),\n    @(\"913\u00f74=\", \"899\u00f78=\"),\n    @(\"628\u00f78=\", \"103\u00f78=\"),\n    @(\"975\u00f79=\", \"983\u00f73=\"),\n    @(\"116\u00f77=\", \"972\u00f77=\"),\n    @(\"386\u00f76=\", \"957\u00f76=\"),\n    @(\"804\u00f78=\", \"893\u00f78=\"),\n    @(\"998\u00f75=\", \"106\u00f75=\"),\n    @(\"404\u00f78=\", \"628\u00f75=\"),\n    @(\"836\u00f72=\", \"863\u00f74=\"),\n    @(\"653\u00f74=\", \"665\u00f76=\"),\n    @(\"637\u00f79=\", \"597\u00f74=\"),\n    @(\"300\u00f76=\", \"763\u00f77=\"),\n    @(\"843\u00f78=\", \"317\u00f73=\"),\n    @(\"430\u00f74=\", \"418\u00f77=\"),\n    @(\"211\u00f73=\", \"585\u00f77=\"),\n    @(\"323\u00f75=\", \"989\u00f78=\"),\n    @(\"676\u00f72=\", \"187\u00f73=\"),\n    @(\"161\u00f74=\", \"999\u00f74=\"),\n    @(\"996\u00f72=\", \"915\u00f75=\"),\n    @(\"412\u00f74=\", \"453\u00f75=\"),\n    @(\"240\u00f72=\", \"497\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $old,      # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $new,      # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
